$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the location names in column A (rows 3-22) from UPPERCASE to Title Case
$ws.Range("A3").Value  = "Alameda"
$ws.Range("A4").Value  = "Antonio Varas"
$ws.Range("A5").Value  = "Campus Arauco"
$ws.Range("A6").Value  = "Campus Villarrica"
$ws.Range("A7").Value  = "Concepción"
$ws.Range("A8").Value  = "Maipú"
$ws.Range("A9").Value  = "Melipilla"
$ws.Range("A10").Value = "Nacimiento"
$ws.Range("A11").Value = "Online"
$ws.Range("A12").Value = "Padre Alonso De Ovalle"
$ws.Range("A13").Value = "Plaza Norte"
$ws.Range("A14").Value = "Plaza Oeste"
$ws.Range("A15").Value = "Plaza Vespucio"
$ws.Range("A16").Value = "Puente Alto"
$ws.Range("A17").Value = "Puerto Montt"
$ws.Range("A18").Value = "San Bernardo"
$ws.Range("A19").Value = "San Carlos De Apoquindo"
$ws.Range("A20").Value = "San Joaquín"
$ws.Range("A21").Value = "Valparaíso"
$ws.Range("A22").Value = "Viña Del Mar"

# Update the view: scroll so row 2 is at the top and select A23
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A23").Select()
